$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Some "Price" column values in this sheet look like plain decimal numbers
# (e.g. "571.94", "7.20", "0.999") even though the source data stores them as
# text strings (other rows use multi-dot "thousands" groupings like
# "65.148.70" which are not valid numbers at all and stay text automatically).
# To keep these particular cells as text too -- matching the original inline
# string storage and avoiding loss of significant trailing/leading zeros -- 
# force a Text number format on them before writing the values, just like you
# would do interactively in Excel.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D11",
    "D12",
    "D13",
    "D14",
    "D17",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D27",
    "D28",
    "D30",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D40",
    "D45",
    "D46",
    "D50",
    "D51",
)
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the cell content/value updates described by the diff.
$ws.Range("D2").Value = '65.148.70'
$ws.Range("E2").Value = '  +1.94%  '
$ws.Range("D3").Value = '3.166.26'
$ws.Range("E3").Value = '  +3.67%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '571.94'
$ws.Range("E5").Value = '  +2.68%  '
$ws.Range("D6").Value = '150.78'
$ws.Range("E6").Value = '  +6.09%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '3.166.01'
$ws.Range("E8").Value = '  +3.73%  '
$ws.Range("E9").Value = '  +2.91%  '
$ws.Range("E10").Value = '  +4.86%  '
$ws.Range("D11").Value = '6.21'
$ws.Range("E11").Value = '  +2.47%  '
$ws.Range("D12").Value = '0.506'
$ws.Range("E12").Value = '  +5.94%  '
$ws.Range("D13").Value = '0.0000274'
$ws.Range("E13").Value = '  +18.81%  '
$ws.Range("D14").Value = '38.19'
$ws.Range("E14").Value = '  +8.61%  '
$ws.Range("D15").Value = '3.682.67'
$ws.Range("E15").Value = '  +3.66%  '
$ws.Range("D16").Value = '65.196.83'
$ws.Range("E16").Value = '  +1.88%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = '7.20'
$ws.Range("E17").Value = '  +6.94%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.168.33'
$ws.Range("E18").Value = '  +3.72%  '
$ws.Range("D20").Value = '511.82'
$ws.Range("E20").Value = '  +7.22%  '
$ws.Range("D21").Value = '14.93'
$ws.Range("E21").Value = '  +7.17%  '
$ws.Range("D22").Value = '15.88'
$ws.Range("E22").Value = '  +10.91%  '
$ws.Range("D23").Value = '0.734'
$ws.Range("E23").Value = '  +8.28%  '
$ws.Range("D24").Value = '7.84'
$ws.Range("E24").Value = '  +3.38%  '
$ws.Range("D25").Value = '84.91'
$ws.Range("E25").Value = '  +3.37%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = '9.14'
$ws.Range("E27").Value = '  +15.15%  '
$ws.Range("D28").Value = '2.90'
$ws.Range("E28").Value = '  +4.03%  '
$ws.Range("E29").Value = '  +8.91%  '
$ws.Range("D30").Value = '27.99'
$ws.Range("E30").Value = '  +6.64%  '
$ws.Range("E31").Value = '  +15.25%  '
$ws.Range("E32").Value = '  +7.96%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").Value = '6.34'
$ws.Range("E34").Value = '  +12.43%  '
$ws.Range("D35").Value = '6.65'
$ws.Range("E35").Value = '  +7.37%  '
$ws.Range("D36").Value = '55.63'
$ws.Range("E36").Value = '  +1.50%  '
$ws.Range("D37").Value = '477.66'
$ws.Range("E37").Value = '  +7.59%  '
$ws.Range("D38").Value = '0.0882'
$ws.Range("E38").Value = '  +9.52%  '
$ws.Range("E39").Value = '  +8.83%  '
$ws.Range("D40").Value = '0.0422'
$ws.Range("E40").Value = '  +3.71%  '
$ws.Range("D41").Value = '3.124.67'
$ws.Range("E41").Value = '  +4.57%  '
$ws.Range("E42").Value = '  +4.97%  '
$ws.Range("E43").Value = '  +6.36%  '
$ws.Range("E44").Value = '  +17.47%  '
$ws.Range("D45").Value = '0.290'
$ws.Range("E45").Value = '  +11.20%  '
$ws.Range("D46").Value = '29.34'
$ws.Range("E46").Value = '  +5.77%  '
$ws.Range("E47").Value = '  +14.32%  '
$ws.Range("E49").Value = '  +2.16%  '
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").Value = '  +11.12%  '
$ws.Range("D51").Value = '123.34'
$ws.Range("E51").Value = '  +4.83%  '
